# Stores maze start and end index to make moving positions more efficient.
# Began implementing an InstructionReceived packet.
#
# Adds two new TODO rows to the "Networks" task table (columns K/L/M) on
# Sheet1, mirroring the existing K/L/M rows above them:
#   Row 18: "Remove rotating box" / "Clean up code used to handle this" / TODO
#   Row 19: Add "CONNECTED" status entry / Or "NOT CONNECTED" / TODO

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 - InstructionReceived packet clean-up task
$ws.Range("K18").Value = "Remove rotating box"
$ws.Range("L18").Value = "Clean up code used to handle this"
$ws.Range("M18").Value = "TODO"
$ws.Range("M18").Interior.Color = $ws.Range("M17").Interior.Color

# Row 19 - connection status entry task
$ws.Range("K19").Value = 'Add "CONNECTED" status entry'
$ws.Range("L19").Value = 'Or "NOT CONNECTED"'
$ws.Range("M19").Value = "TODO"
$ws.Range("M19").Interior.Color = $ws.Range("M17").Interior.Color

# Match the author's final selection in the saved workbook
$ws.Range("M19").Select() | Out-Null
